$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the empty/duplicate trailing columns (H, I, J) ---
$ws.Range("H1:J4").EntireColumn.Delete()

# --- Header row (plain text, safe to assign directly) ---
$ws.Range("B1").Value = "Academic research & education"
$ws.Range("C1").Value = "Social causes"
$ws.Range("D1").Value = "Arts"
$ws.Range("E1").Value = "Politics"
$ws.Range("F1").Value = "Sport"
$ws.Range("G1").Value = "Total amount donated (€)"

# --- Numeric-looking values must stay text, so stage them via formulas in a
#     scratch area (laid out row-by-row like the destination), then
#     copy/paste-special so the pasted cells keep the "text" cell type
#     instead of Excel auto-converting them to numbers ---
$ws.Range("Z2").Formula  = '="59.5"'
$ws.Range("AA2").Formula = '="17.7"'
$ws.Range("AB2").Formula = '="12.5"'
$ws.Range("AC2").Formula = '="6.3"'
$ws.Range("AD2").Formula = '="4.1"'
$ws.Range("AE2").Formula = '="318,376"'

$ws.Range("Z3").Formula  = '="68.4"'
$ws.Range("AA3").Formula = '="18.6"'
$ws.Range("AB3").Formula = '="5.7"'
$ws.Range("AC3").Formula = '="5.3"'
$ws.Range("AD3").Formula = '="1.9"'
$ws.Range("AE3").Formula = '="656,851"'

$ws.Range("Z4").Formula  = '="43.8"'
$ws.Range("AA4").Formula = '="21.7"'
$ws.Range("AB4").Formula = '="25.2"'
$ws.Range("AC4").Formula = '="6.3"'
$ws.Range("AD4").Formula = '="3.0"'
$ws.Range("AE4").Formula = '="657,304"'

$ws.Range("Z2:AE4").Copy()
$ws.Range("B2").PasteSpecial()

$ws.Range("Z2:AE4").ClearContents()
